# Word COM-interop script implementing the Black-Holes -> Chemistry rewrite
$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $old"
    }
    $target = $d.Range($rng.Start, $rng.End)
    $target.Text = $new
    return $target
}

function Insert-After-Text($afterRange, $text) {
    $ins = $d.Range($afterRange.End, $afterRange.End)
    $ins.InsertAfter($text)
    $newEnd = $afterRange.End + $text.Length
    return $d.Range($afterRange.End, $newEnd)
}

# --- Title / byline / email ---
Replace-Exact "The Enigma of Black Holes: Unveiling Cosmic Mysteries" "Unraveling the Composition of Matter: A Dive Into Chemistry" | Out-Null
Replace-Exact "Alex Richards" "Marcus King" | Out-Null
Replace-Exact "alexrichards@astronews" "marcusking@gmail" | Out-Null
Replace-Exact "net" "com" | Out-Null

# --- Body paragraph 1 ---
Replace-Exact "In the vast expanse of the universe, there lurk enigmatic entities known as black holes, captivating the scientific world with their profound mysteries" "Chemistry, the study of substances, their properties, and how they change, unveils a fascinating realm of interactions at the molecular level" | Out-Null

Replace-Exact " These celestial behemoths, born from the gravitational collapse of massive stars, defy our current understanding of physics and challenge our perception of reality" " As we delve into the intricacies of matter, we are awestruck by its fundamental building blocks, the elements, and the countless compounds they forge" | Out-Null

$r9 = Replace-Exact " The quest to unveil the secrets of black holes has spurred an intense scientific endeavor, as we strive to comprehend their perplexing properties and their profound implications for the cosmos" " Chemistry orchestrates the symphony of life, governing chemical reactions essential for biological processes"
$r9b = Insert-After-Text $r9 "."
Insert-After-Text $r9b " From the respiration that fuels our bodies to the food we consume, chemistry plays an indispensable role in the survival of all living organisms" | Out-Null

Replace-Exact "From the dawn of their theoretical conception, black holes have sparked fervent debate and speculation among physicists, astronomers, and cosmologists" "Further exploration reveals the composition of matter, unveiling the subatomic universe" | Out-Null

Replace-Exact " Their gravitational pull is so intense that not even light can escape their clutches, creating an event horizon that marks the boundary of no return" " Electrons, protons, and neutrons dance within atoms, defining their properties and behavior" | Out-Null

$r15 = Replace-Exact " Within this enigmatic region, the laws of physics, as we know them, break down, giving rise to a realm of uncertainty and theoretical paradoxes" " The periodic table, a roadmap of the elements, categorizes these fundamental components based on their properties, allowing us to predict their reactivity and tendencies"
$r15b = Insert-After-Text $r15 "."
Insert-After-Text $r15b " As we unravel the mysteries of atomic structure, we gain insight into the forces that govern chemical bonding, shaping the diverse materials that make up our world" | Out-Null

Replace-Exact "The study of black holes has opened new avenues of exploration, pushing the boundaries of our scientific understanding" "Through experimentation and observation, chemists unlock the secrets of chemical reactions" | Out-Null

Replace-Exact " By observing and analyzing these cosmic phenomena, we have gained insights into the behavior of matter under extreme conditions, the nature of spacetime, and the fundamental forces that govern the universe" " They investigate the conditions that initiate and control these transformations, manipulating temperature, pressure, and catalysts to direct the course of reactions" | Out-Null

$r21 = Replace-Exact " The pursuit of unraveling the mysteries of black holes continues to ignite the imagination and fuel the progress of science, promising groundbreaking discoveries and a deeper understanding of our place in the cosmos" " Understanding chemical reactions empowers us to harness their potential, creating new substances, developing pharmaceuticals, and addressing environmental challenges"
$r21b = Insert-After-Text $r21 "."
Insert-After-Text $r21b " By understanding the language of chemistry, we unlock the door to countless technological advancements that enhance our lives" | Out-Null

# --- Summary heading stays the same ---

# --- Summary paragraph ---
Replace-Exact "Black holes, with their enigmatic nature and profound implications, have captured the scientific community's attention" "Chemistry, in its essence, uncovers the secrets of matter and its transformations" | Out-Null

Replace-Exact " These cosmic entities challenge our understanding of physics and offer a glimpse into the extremes of the universe" " It delves into the composition of substances, revealing the fundamental building blocks of the universe, the elements, and their intricate interactions" | Out-Null

Replace-Exact " Through ongoing research " " By dissecting chemical reactions " | Out-Null

$r29 = Replace-Exact "and observation, we continue to uncover the mysteries surrounding black holes, shedding light on their behavior, properties, and the fundamental forces that shape our cosmos" "and manipulating their conditions, chemists unlock the potential for creating new materials, advancing medical treatments, and addressing pressing environmental issues"
$r29b = Insert-After-Text $r29 "."
Insert-After-Text $r29b " Chemistry serves as a gateway to understanding the fabric of our world, empowering us to harness the power of matter and shape a better future" | Out-Null

# --- Add trailing empty paragraph at the very end of the document body ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Host "Done"
